# The document originally started with a long block of draft/working
# paragraphs (requirements notes, a page break, etc.) before the real
# "Contribuciones" heading. This edit strips all of that preliminary
# content so the document begins directly at the "Contribuciones"
# section heading (turning the file into the clean master document).

$d = $word.ActiveDocument

# Locate the "Contribuciones" heading paragraph by searching for its
# text rather than hard-coding a paragraph index, so the script is
# resilient to how the engine enumerates paragraphs.
$headingRange = $d.Content.Duplicate
$found = $headingRange.Find.Execute("Contribuciones", $true, $true, $false,
                                     $false, $false, $true, 1, $false,
                                     "", 0)

if ($found -and $headingRange.Start -gt 0) {
    # Remove everything from the very start of the document up to (but
    # not including) the heading paragraph.
    $deleteRange = $d.Range(0, $headingRange.Start)
    $deleteRange.Delete()
}
